$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 31 de Agosto de 2020 a las 04:36"

# Bolivia (row 30) - updated case counts
$ws.Range("B30").Value = 115968
$ws.Range("C30").Value = 614
$ws.Range("D30").Value = 58497
$ws.Range("E30").Value = 52505
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 28
$ws.Range("H30").Value = 4966

# Honduras (row 50) - updated case counts
$ws.Range("B50").Value = 60174
$ws.Range("C50").Value = 529
$ws.Range("D50").Value = 10242
$ws.Range("E50").Value = 48074
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 16
$ws.Range("H50").Value = 1858

# Australia (row 72) - updated case counts
$ws.Range("B72").Value = 25746
$ws.Range("C72").Value = 76
$ws.Range("D72").Value = 21116
$ws.Range("E72").Value = 3978
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 41
$ws.Range("H72").Value = 652

# Row 76 now holds "Corea del Sur" (its case count overtook Bosnia's, so it
# moves up in the descending sort) with its updated figures.
$ws.Range("A76").Value = "Corea del Sur"
$ws.Range("B76").Value = 19947
$ws.Range("C76").Value = 248
$ws.Range("D76").Value = 14973
$ws.Range("E76").Value = 4650
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 324

# Row 77 now holds "Bosnia y Herzegovina" with its (unchanged) figures.
$ws.Range("A77").Value = "Bosnia y Herzegovina"
$ws.Range("B77").Value = 19793
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 12961
$ws.Range("E77").Value = 6234
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 598

# Belice (row 165) - updated case counts
$ws.Range("B165").Value = 993
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 117
$ws.Range("E165").Value = 863
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 13
